{"js": "// Change title to include \"Farm Owner and\" before \"Farm Manager Version)\"\nconst body = context.document.body;\nconst results = body.search(\"Informed Consent Form and Privacy Policy (Farm Manager Version)\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const range = results.items[0];\n  // Insert the replacement text, preserving the original formatting of the run.\n  range.insertText(\n    \"Informed Consent Form and Privacy Policy (Farm Owner and Farm Manager Version)\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Change title to include \"Farm Owner and\" before \"Farm Manager Version)\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"Informed Consent Form and Privacy Policy (Farm Manager Version)\",\n    $false, $false, $false, $false, $false, $true, 0, $false,\n    \"Informed Consent Form and Privacy Policy (Farm Owner and Farm Manager Version)\",\n    2\n) | Out-Null\n"}
